$d = $word.ActiveDocument

# 1) Rename bookmarks (structural, not text content)
$d.Bookmarks.Item("masterminter-contract").Name = "draft-master-minter-contract"
$d.Bookmarks.Item("interaction-with-natgxtoken-contract").Name = "interaction-with-natgx-token-contract"
$d.Bookmarks.Item("configuring-the-masterminter").Name = "configuring-the-master-minter"
$d.Bookmarks.Item("masterminter-vs.-mintcontroller").Name = "master-minter-vs.-mintcontroller"

# 2) Fix the title heading (first paragraph only): "MasterMinter contract" -> "DRAFT: Master Minter contract"
$headingRange = $d.Paragraphs(1).Range
$headingRange.Find.Execute("MasterMinter contract", $true, $true, $false, $false, $false, $true, 1, $false, "DRAFT: Master Minter contract", 2)

# 3) Globally insert a space between "Master" and "Minter" wherever "MasterMinter" appears
#    (case-sensitive so the lower-case role name "masterMinter" is left untouched).
$d.Content.Find.Execute("MasterMinter", $true, $true, $false, $false, $false, $true, 1, $false, "Master Minter", 2)

# 4) Globally insert a space between "NATGX" and "Token" wherever "NATGXToken" appears.
$d.Content.Find.Execute("NATGXToken", $true, $true, $false, $false, $false, $true, 1, $false, "NATGX Token", 2)
